# ============================================================================
# znfmd630 - TMS Ordens de Frete: documentation update
# ============================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$wsEnums = $wb.Worksheets.Item("enums")

# ----------------------------------------------------------------------------
# 1) Translate / rename column headers and enum-domain labels (PT-BR pass)
# ----------------------------------------------------------------------------
$ws.Range("J1").Value = "Nº seq"
$ws.Range("L1").Value = "Campo combinado"
$ws.Range("P1").Value = "Comp"

$ws.Range("O4").Value  = "Inteiro"
$ws.Range("O46").Value = "Inteiro"
$ws.Range("O49").Value = "Inteiro"

$ws.Range("O27").Value = "Enumerado"
$ws.Range("O43").Value = "Enumerado"
$ws.Range("O54").Value = "Enumerado"
$ws.Range("O56").Value = "Enumerado"

$ws.Range("O57").Value = "Data"

$ws.Range("O15").Value = "Data/Hora UTC"
$ws.Range("O30").Value = "Data/Hora UTC"
$ws.Range("O52").Value = "Data/Hora UTC"

$wsEnums.Range("D4").Value  = "Inteiro"
$wsEnums.Range("D9").Value  = "Enumerado"
$wsEnums.Range("D10").Value = "Conjunto"
$wsEnums.Range("D11").Value = "Data"
$wsEnums.Range("D12").Value = "Texto"
$wsEnums.Range("D14").Value = "Data/Hora UTC"
$wsEnums.Range("D16").Value = "Objeto binário grande"

# ----------------------------------------------------------------------------
# 2) Row 28 (stat.c) re-pointed to a new domain + label updates
# ----------------------------------------------------------------------------
$ws.Range("M28").Value = "zn"
$ws.Range("N28").Value = "fmd.stof.c"
$ws.Range("O28").Value = "Enumerado"
$ws.Range("P28").Value = 7
$ws.Range("Q28").Value = "znfmd200.stat.c"
$ws.Range("R28").Value = "Status"

# ----------------------------------------------------------------------------
# 3) Row 42 (fdtc.c) label updates
# ----------------------------------------------------------------------------
$ws.Range("Q42").Value = "sls.obsoleto.c"
$ws.Range("R42").Value = "Obsoleto"

# ----------------------------------------------------------------------------
# 4) Append 12 new field definitions (rows 59-70), cloning row 58's
#    cell styles/types via copy, then patching per-row values.
# ----------------------------------------------------------------------------
function Add-FieldRow {
    param(
        [int]$RowNum,
        [int]$Seq,
        [string]$Field,
        [string]$Domain,
        [string]$DataType,
        [int]$Length,
        [string]$Label,
        [string]$Description
    )

    $srcRow = $RowNum - 1
    $ws.Range("A${srcRow}:R${srcRow}").Copy()
    $ws.Range("A${RowNum}:R${RowNum}").PasteSpecial(-4104)
    $excel.CutCopyMode = $false

    $ws.Range("J${RowNum}").Value = $Seq
    $ws.Range("K${RowNum}").Value = $Field
    $ws.Range("N${RowNum}").Value = $Domain
    $ws.Range("O${RowNum}").Value = $DataType
    $ws.Range("P${RowNum}").Value = $Length
    $ws.Range("Q${RowNum}").Value = $Label
    $ws.Range("R${RowNum}").Value = $Description
}

Add-FieldRow 59 58 "vllq.c" "amnt"     "Double"       19 "sls.obsoleto.c"  "Obsoleto"
Add-FieldRow 60 59 "copo.c" "mcs.int2" "Inteiro"       2 "znsls401.copo.c" "Coleta ou Postagem"
Add-FieldRow 61 60 "frpe.c" "amnt"     "Double"       19 "znfmd630.frpe.c" "Frete Peso"
Add-FieldRow 62 61 "dtco.c" "date"     "Data/Hora UTC" 22 "znfmd630.dtco.c" "Data Corrigida"
Add-FieldRow 63 62 "advc.c" "amnt"     "Double"       19 "znfmd630.advc.c" "Ad Valorem Calculado"
Add-FieldRow 64 63 "pedc.c" "amnt"     "Double"       19 "znfmd630.pedc.c" "Pedágio Calculado"
Add-FieldRow 65 64 "reet.c" "yesno"    "Enumerado"     5 "znfmd630.reet.c" "Reentrega Gerada"
Add-FieldRow 66 65 "sqrt.c" "mcs.long" "Long"         10 "znfmd630.sqrt.c" "Sequencial Reentrega"
Add-FieldRow 67 66 "rcal.c" "yesno"    "Enumerado"     5 "znfmd630.rcal.c" "Recalculado"
Add-FieldRow 68 67 "dtrc.c" "date"     "Data/Hora UTC" 22 "znfmd630.dtrc.c" "Data do Recalculo"
Add-FieldRow 69 68 "udap.c" "date"     "Data/Hora UTC" 22 "znfmd630.udap.c" "Última Data Acompanhar Pedido"
Add-FieldRow 70 69 "dtpe.c" "date"     "Data/Hora UTC" 22 "znfmd630.dtpe.c" "Data Prevista de entrega"

# ----------------------------------------------------------------------------
# 5) Refresh the O2:O1048575 list validation so it also covers the new rows,
#    and move the selection to K2 as in the edited file.
# ----------------------------------------------------------------------------
$ws.Range("O2:O1048575").Validation.Delete()
$ws.Range("O2:O1048575").Validation.Add(3, 1, 1, "ttadv.type")

$ws.Range("K2").Select()

# ----------------------------------------------------------------------------
# 6) Re-fit column widths for the columns whose content changed width.
# ----------------------------------------------------------------------------
$ws.Columns.Item(10).AutoFit() | Out-Null
$ws.Columns.Item(12).AutoFit() | Out-Null
$ws.Columns.Item(15).AutoFit() | Out-Null
$ws.Columns.Item(16).AutoFit() | Out-Null
